$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert three brand-new earthquake rows right under the header row,
#     shifting every existing data row down by three. Excel's default
#     Insert shifts the existing rows down (xlShiftDown). ---
$ws.Rows("2:4").Insert()

# Scratch cell used to "type" a value through a formula (so the engine's
# literal-entry type-inference - which would otherwise convert a numeric-
# looking string like "15.00" into the number 15 and bump the cell style -
# never runs). Copy/PasteSpecial(values) transfers the already-computed
# text result verbatim, keeping the cell General-formatted/unstyled, just
# like every other text cell in this sheet.
$scratch = $ws.Cells.Item(1, 20)

function Set-TextValue($row, $col, $text) {
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4163)  # xlPasteValues
}

# Row 2 - new entry: M5.0, 越南 (Vietnam)
Set-TextValue 2 1 "5.0"
Set-TextValue 2 2 "2024-07-28 12:35:13"
Set-TextValue 2 3 "15.00"
Set-TextValue 2 4 "108.20"
Set-TextValue 2 5 "10"
Set-TextValue 2 6 "越南"

# Row 3 - new entry: M3.4, 新疆巴音郭楞州尉犁县
Set-TextValue 3 1 "3.4"
Set-TextValue 3 2 "2024-07-28 05:41:53"
Set-TextValue 3 3 "40.84"
Set-TextValue 3 4 "84.10"
Set-TextValue 3 5 "21"
Set-TextValue 3 6 "新疆巴音郭楞州尉犁县"

# Row 4 - new entry: M3.1, 西藏那曲市双湖县
Set-TextValue 4 1 "3.1"
Set-TextValue 4 2 "2024-07-28 05:22:21"
Set-TextValue 4 3 "33.34"
Set-TextValue 4 4 "87.27"
Set-TextValue 4 5 "10"
Set-TextValue 4 6 "西藏那曲市双湖县"

# Clean up the scratch cell so it leaves no trace in the saved workbook.
$scratch.ClearContents()

# The sheet's used range must stay fixed at A1:F82 (oldest rows fall off
# the bottom of the feed), so remove the three rows pushed past row 82.
$ws.Rows("83:85").Delete()
